$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns: AD = Wins, AE = Losses, AF = Ties
# Header row (row 1) - copy formatting from the existing header cell (AC1)
# so the new headers match the bold/bordered/centered header style.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-40: every row gets the same season record (83 wins, 79 losses, 0 ties)
for ($row = 2; $row -le 40; $row++) {
    $ws.Cells.Item($row, 30).Value = 83   # AD
    $ws.Cells.Item($row, 31).Value = 79   # AE
    $ws.Cells.Item($row, 32).Value = 0    # AF
}
